# Auto-generated edit script: updates market-price / profit columns (H:N)
# across the Hyperion_Profits worksheets, per the scheduled-runner refresh diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 645
$ws.Range("I19").Value = 98.8
$ws.Range("J19").Value = 1035.1428
$ws.Range("K19").Value = 98.8
$ws.Range("L19").Value = 1035.1428
$ws.Range("M19").Value = 76.2
$ws.Range("N19").Value = -1385.1428
$ws.Range("H40").Value = 3117.4285
$ws.Range("I40").Value = 3296.7646
$ws.Range("J40").Value = 2355.25
$ws.Range("K40").Value = 3296.7646
$ws.Range("L40").Value = 2355.25
$ws.Range("M40").Value = -3121.7646
$ws.Range("N40").Value = -2705.25
$ws.Range("H86").Value = 3812.12
$ws.Range("I86").Value = 4263.8
$ws.Range("K86").Value = 4263.8
$ws.Range("M86").Value = -3140.8
$ws.Range("H89").Value = 3812.12
$ws.Range("I89").Value = 4263.8
$ws.Range("K89").Value = 21319
$ws.Range("M89").Value = -15703
$ws.Range("H99").Value = 380.6
$ws.Range("J99").Value = 514.2
$ws.Range("L99").Value = 1542.6
$ws.Range("N99").Value = -4538.6
$ws.Range("H106").Value = 1827.4445
$ws.Range("I106").Value = 1485.5714
$ws.Range("K106").Value = 1485.5714
$ws.Range("M106").Value = -854.5714
$ws.Range("H113").Value = 5122.64
$ws.Range("I113").Value = 4859.4546
$ws.Range("J113").Value = 5329.4287
$ws.Range("K113").Value = 4859.4546
$ws.Range("L113").Value = 5329.4287
$ws.Range("M113").Value = -1605.4546
$ws.Range("N113").Value = -11837.4287
$ws.Range("H138").Value = 2664.8281
$ws.Range("I138").Value = 917.3043
$ws.Range("J138").Value = 3645.1462
$ws.Range("K138").Value = 2751.9129
$ws.Range("L138").Value = 10935.4386
$ws.Range("M138").Value = 2388.0871
$ws.Range("N138").Value = -21215.4386
$ws.Range("H141").Value = 10118.6
$ws.Range("I141").Value = 10484.214
$ws.Range("K141").Value = 31452.642
$ws.Range("M141").Value = -26272.642

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 36375
$ws.Range("J51").Value = 36375
$ws.Range("L51").Value = 36375
$ws.Range("N51").Value = -37887
$ws.Range("H97").Value = 1618362.9
$ws.Range("I97").Value = 2941753
$ws.Range("J97").Value = 886.2222
$ws.Range("K97").Value = 2941753
$ws.Range("L97").Value = 886.2222
$ws.Range("M97").Value = -2941257
$ws.Range("N97").Value = -1878.2222
$ws.Range("H102").Value = 3791643
$ws.Range("I102").Value = 4905189
$ws.Range("K102").Value = 4905189
$ws.Range("M102").Value = -4903567

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19610792
$ws.Range("I20").Value = 47621868
$ws.Range("K20").Value = 47621868
$ws.Range("M20").Value = -47621621
$ws.Range("H22").Value = 3362.125
$ws.Range("I22").Value = 4279.4
$ws.Range("K22").Value = 4279.4
$ws.Range("M22").Value = -4106.4
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H59").Value = 122784.664
$ws.Range("J59").Value = 122784.664
$ws.Range("L59").Value = 122784.664
$ws.Range("N59").Value = -124478.664
$ws.Range("H80").Value = 454.08572
$ws.Range("J80").Value = 444.38095
$ws.Range("L80").Value = 444.38095
$ws.Range("N80").Value = -2440.38095
$ws.Range("H83").Value = 454.08572
$ws.Range("J83").Value = 444.38095
$ws.Range("L83").Value = 2221.90475
$ws.Range("N83").Value = -12205.90475
$ws.Range("H99").Value = 7146560.5
$ws.Range("I99").Value = 9527248
$ws.Range("J99").Value = 4499.4
$ws.Range("K99").Value = 9527248
$ws.Range("L99").Value = 4499.4
$ws.Range("M99").Value = -9525750
$ws.Range("N99").Value = -7495.4
$ws.Range("H134").Value = 5073.3887
$ws.Range("I134").Value = 4497.4
$ws.Range("K134").Value = 13492.2
$ws.Range("M134").Value = -10957.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24710.268
$ws.Range("I31").Value = 2925.3333
$ws.Range("J31").Value = 30156.5
$ws.Range("K31").Value = 2925.3333
$ws.Range("L31").Value = 30156.5
$ws.Range("M31").Value = -2630.3333
$ws.Range("N31").Value = -30746.5
$ws.Range("H34").Value = 24710.268
$ws.Range("I34").Value = 2925.3333
$ws.Range("J34").Value = 30156.5
$ws.Range("K34").Value = 2925.3333
$ws.Range("L34").Value = 30156.5
$ws.Range("M34").Value = -2723.3333
$ws.Range("N34").Value = -30560.5
$ws.Range("H134").Value = 3401.3125
$ws.Range("I134").Value = 2306.9524
$ws.Range("K134").Value = 6920.8572
$ws.Range("M134").Value = -4385.8572

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 209.78
$ws.Range("I2").Value = 102.90625
$ws.Range("K2").Value = 617.4375
$ws.Range("M2").Value = -504.4375
$ws.Range("H3").Value = 1609.25
$ws.Range("I3").Value = 848.2222
$ws.Range("J3").Value = 3892.3333
$ws.Range("K3").Value = 2544.6666
$ws.Range("L3").Value = 11676.9999
$ws.Range("M3").Value = -2432.6666
$ws.Range("N3").Value = -11900.9999
$ws.Range("H5").Value = 72466.64
$ws.Range("I5").Value = 1015.2222
$ws.Range("K5").Value = 3045.6666
$ws.Range("M5").Value = -2933.6666
$ws.Range("H7").Value = 2670.625
$ws.Range("J7").Value = 1687.3334
$ws.Range("L7").Value = 5062.0002
$ws.Range("N7").Value = -5286.0002
$ws.Range("H46").Value = 198513.3
$ws.Range("I46").Value = 1111277.6
$ws.Range("J46").Value = 2920.9285
$ws.Range("K46").Value = 3333832.8
$ws.Range("L46").Value = 8762.7855
$ws.Range("M46").Value = -3333741.8
$ws.Range("N46").Value = -8944.7855
$ws.Range("H63").Value = 11502.375
$ws.Range("I63").Value = 3019
$ws.Range("J63").Value = 12714.286
$ws.Range("K63").Value = 9057
$ws.Range("L63").Value = 38142.858
$ws.Range("M63").Value = -8308
$ws.Range("N63").Value = -39640.858
$ws.Range("H66").Value = 11502.375
$ws.Range("I66").Value = 3019
$ws.Range("J66").Value = 12714.286
$ws.Range("K66").Value = 27171
$ws.Range("L66").Value = 114428.574
$ws.Range("M66").Value = -23427
$ws.Range("N66").Value = -121916.574
$ws.Range("H86").Value = 55.5
$ws.Range("I86").Value = 55.5
$ws.Range("K86").Value = 166.5
$ws.Range("M86").Value = 1019.5
$ws.Range("H89").Value = 55.5
$ws.Range("I89").Value = 55.5
$ws.Range("K89").Value = 499.5
$ws.Range("M89").Value = 5428.5
$ws.Range("H92").Value = 908.8
$ws.Range("I92").Value = 274.75
$ws.Range("J92").Value = 1331.5
$ws.Range("K92").Value = 824.25
$ws.Range("L92").Value = 3994.5
$ws.Range("M92").Value = 423.75
$ws.Range("N92").Value = -6490.5
$ws.Range("H135").Value = 72466.64
$ws.Range("I135").Value = 1015.2222
$ws.Range("K135").Value = 9136.9998
$ws.Range("M135").Value = -6601.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 33356666
$ws.Range("I49").Value = 100000000
$ws.Range("K49").Value = 100000000
$ws.Range("M49").Value = -99999816

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 69352.08
$ws.Range("I22").Value = 111922.5
$ws.Range("J22").Value = 1239.4
$ws.Range("K22").Value = 111922.5
$ws.Range("L22").Value = 1239.4
$ws.Range("M22").Value = -111627.5
$ws.Range("N22").Value = -1829.4
$ws.Range("H27").Value = 69352.08
$ws.Range("I27").Value = 111922.5
$ws.Range("J27").Value = 1239.4
$ws.Range("K27").Value = 111922.5
$ws.Range("L27").Value = 1239.4
$ws.Range("M27").Value = -111815.5
$ws.Range("N27").Value = -1453.4
$ws.Range("H46").Value = 7189.7
$ws.Range("J46").Value = 9142.429
$ws.Range("L46").Value = 9142.429
$ws.Range("N46").Value = -9518.429
$ws.Range("H54").Value = 31500
$ws.Range("J54").Value = 31500
$ws.Range("L54").Value = 31500
$ws.Range("N54").Value = -32788
$ws.Range("H132").Value = 18497.75
$ws.Range("I132").Value = 19711.715
$ws.Range("K132").Value = 59135.145
$ws.Range("M132").Value = -56605.145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 20145.285
$ws.Range("J31").Value = 21254.5
$ws.Range("L31").Value = 21254.5
$ws.Range("N31").Value = -21950.5
$ws.Range("H81").Value = 7578903.5
$ws.Range("I81").Value = 13890433
$ws.Range("J81").Value = 5068
$ws.Range("K81").Value = 27780866
$ws.Range("L81").Value = 10136
$ws.Range("M81").Value = -27779805
$ws.Range("N81").Value = -12258
$ws.Range("H84").Value = 7578903.5
$ws.Range("I84").Value = 13890433
$ws.Range("J84").Value = 5068
$ws.Range("K84").Value = 138904330
$ws.Range("L84").Value = 50680
$ws.Range("M84").Value = -138899026
$ws.Range("N84").Value = -61288
$ws.Range("H96").Value = 3206.8262
$ws.Range("I96").Value = 2938.8823
$ws.Range("J96").Value = 3966
$ws.Range("K96").Value = 2938.8823
$ws.Range("L96").Value = 3966
$ws.Range("M96").Value = -1565.8823
$ws.Range("N96").Value = -6712
$ws.Range("H100").Value = 2782.25
$ws.Range("I100").Value = 3932
$ws.Range("J100").Value = 866
$ws.Range("K100").Value = 7864
$ws.Range("L100").Value = 1732
$ws.Range("M100").Value = -7323
$ws.Range("N100").Value = -2814
$ws.Range("H132").Value = 12799084
$ws.Range("I132").Value = 15628372
$ws.Range("J132").Value = 727453.4399999999
$ws.Range("K132").Value = 46885116
$ws.Range("L132").Value = 2182360.32
$ws.Range("M132").Value = -46882586
$ws.Range("N132").Value = -2187420.32

